# Generate inputs for 2 scenarios to test importance of ET daily timestep.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two date values (B3, B4) from 8/15/2019 (43692) to 8/19/2019 (43696)
$ws.Range("B3").Value = 43696
$ws.Range("B4").Value = 43696

# Set column A width (closest value the host's character-width
# quantization can reach to the target stored width of 15.5546875)
$ws.Columns.Item(1).ColumnWidth = 14.6666666666667

# Update the active selection
$ws.Range("C9").Select()
